$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.661.28'
$ws.Range("E2").Value = '  -1.22%  '

$ws.Range("D3").Value = '1.860.96'
$ws.Range("E3").Value = '  -1.91%  '

$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '334.59'
$ws.Range("E5").Value = '  -0.48%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.011'
$ws.Range("E6").Value = '  -0.32%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4626'
$ws.Range("E7").Value = '  -1.82%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3910'
$ws.Range("E8").Value = '  -0.84%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '46.39'
$ws.Range("E9").Value = '  -2.44%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07955'
$ws.Range("E10").Value = '  -1.66%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.000'
$ws.Range("E11").Value = '  -2.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '21.64'
$ws.Range("E12").Value = '  -1.51%  '

$ws.Range("D13").Value = '1.868.60'
$ws.Range("E13").Value = '  -0.65%  '

$ws.Range("E14").Value = '  -0.69%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.212'
$ws.Range("E15").Value = '  +0.61%  '

$ws.Range("E16").Value = '  -0.25%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '88.29'
$ws.Range("E17").Value = '  +0.88%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06724'
$ws.Range("E18").Value = '  -1.11%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.00001041'
$ws.Range("E19").Value = '  -1.05%  '

$ws.Range("E20").Value = '  -0.51%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.013'
$ws.Range("E21").Value = '  -0.18%  '

$ws.Range("D22").Value = '27.671.62'
$ws.Range("E22").Value = '  -1.14%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.441'
$ws.Range("E23").Value = '  -1.63%  '

$ws.Range("E24").Value = '  -1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.311'
$ws.Range("E25").Value = '  -1.54%  '

$ws.Range("D26").Value = '2.082.69'
$ws.Range("E26").Value = '  -1.05%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '159.66'
$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  -2.28%  '

$ws.Range("E29").Value = '  +1.62%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.445'
$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '121.93'
$ws.Range("E31").Value = '  -0.24%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9764'
$ws.Range("E32").Value = '  -0.10%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09395'
$ws.Range("E33").Value = '  -1.41%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.626'
$ws.Range("E34").Value = '  -0.56%  '

$ws.Range("E35").Value = '  -1.62%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.330'
$ws.Range("E36").Value = '  -5.91%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02236'
$ws.Range("E37").Value = '  -1.43%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06012'
$ws.Range("E38").Value = '  -2.37%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '8.327'
$ws.Range("E39").Value = '  +2.79%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.194'
$ws.Range("E40").Value = '  -2.61%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.011'
$ws.Range("E41").Value = '  -0.27%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5942'
$ws.Range("E42").Value = '  -1.26%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1868'
$ws.Range("E43").Value = '  -1.32%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.34'
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.247'
$ws.Range("E45").Value = '  -1.52%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.5601'
$ws.Range("E46").Value = '  -2.17%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '12.09'
$ws.Range("E47").Value = '  -0.88%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.917'
$ws.Range("E48").Value = '  -1.54%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.06730'
$ws.Range("E49").Value = '  -3.05%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '111.64'
$ws.Range("E50").Value = '  -2.36%  '

$ws.Range("E51").Value = '  -2.01%  '
